# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-31 22:15:36
#
# The "Recorded By" column (column G) lists the users/systems that touched a
# given attendance record as a comma-separated string. Upstream normalised
# the ordering of that list so the first contributor moves to the end
# (a left-rotation by one element) on every data row that has more than one
# name. Single-name rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the "Recorded By" column from the header row (defaults to column 7 /
# "G", which is where it lives in this report, but search to stay robust).
$headerRow = 1
$recordedByCol = 7
$usedCols = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $usedCols; $c++) {
    $headerVal = $ws.Cells.Item($headerRow, $c).Value2
    if ($headerVal -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

# Find the last used row so we cover every data row regardless of sheet size.
$lastRow = $ws.UsedRange.Rows.Count
if ($ws.UsedRange.Row -gt 1) {
    $lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
}

for ($r = ($headerRow + 1); $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
            $cell.Value = $rotated
        }
    }
}
